# "Rifatti test 1 2 3 4 5" - redo test cases 1-5 (rows 8-12) with fresh
# execution timestamps / traceId / workflowInstanceId, all re-run on
# 2023-03-28, and move the active selection to I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 8  - Test case 1 (VALIDAZIONE_CDA2_LAB_CT1)
$ws.Range("F8").Value = "2023-03-28"
$ws.Range("G8").Value = "2023-03-28T08:13:57Z"
$ws.Range("H8").Value = "b90aa6c5e137fb10"
$ws.Range("I8").Value = "2.16.840.1.113883.2.9.2.30.1c1a7089bb719a940d221bde08f7b44fac0fda9daa2d14ec5030d1de569181cc.f87344359c"

# Row 9  - Test case 2 (VALIDAZIONE_CDA2_LAB_CT2)
$ws.Range("F9").Value = "2023-03-28"
$ws.Range("G9").Value = "2023-03-28T08:15:39Z"
$ws.Range("H9").Value = "40eef6f3fec6cff0"
$ws.Range("I9").Value = "2.16.840.1.113883.2.9.2.30.da9c817519f1a1d3fa08df87ef4f60470a1cd221030d816d4ce28995ea526188.c835ffde29"

# Row 10 - Test case 3 (VALIDAZIONE_CDA2_LAB_CT3)
$ws.Range("F10").Value = "2023-03-28"
$ws.Range("G10").Value = "2023-03-28T08:16:51Z"
$ws.Range("H10").Value = "6a13418da786185b"
$ws.Range("I10").Value = "2.16.840.1.113883.2.9.2.30.d60f69eae987ffad24e25fe40cd4d7b50b1f93aff9845aecf91e84d67983a46e.0cc0cafd68"

# Row 11 - Test case 4 (VALIDAZIONE_CDA2_LAB_CT4)
$ws.Range("F11").Value = "2023-03-28"
$ws.Range("G11").Value = "2023-03-28T08:18:03Z"
$ws.Range("H11").Value = "6ef2b0d6e831b598"
$ws.Range("I11").Value = "2.16.840.1.113883.2.9.2.30.e58c5be5971a857a586bd84787f721d7d4df18197996c7307351686e34f63e12.551e568be8"

# Row 12 - Test case 5 (VALIDAZIONE_CDA2_LAB_CT5)
$ws.Range("F12").Value = "2023-03-28"
$ws.Range("G12").Value = "2023-03-28T08:19:16Z"
$ws.Range("H12").Value = "fc63ea260247fc41"
$ws.Range("I12").Value = "2.16.840.1.113883.2.9.2.30.7506cfa3ef9712317a27b582e4b35f11c815c6745d9fdbb0aa668cf8b99a860c.d74ceb3c6e"

# Move/collapse the saved selection to I12 (also clears the scrolled
# topLeftCell from the prior view).
$ws.Range("I12").Select()
